$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) values of the columns that move together
# for rows 17..22 (A, B, D, E, F, G, H, P, Q, R, S, Z, AB)
$cols = @("A","B","D","E","F","G","H","P","Q","R","S","Z","AB")
$rows = 17..22

$data = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $data[$r] = $rowData
}

# Mapping: new row -> old row (where the content now living at new row came from)
$map = @{
    17 = 22
    18 = 17
    19 = 21
    20 = 19
    21 = 20
    22 = 18
}

foreach ($newRow in $rows) {
    $oldRow = $map[$newRow]
    $src = $data[$oldRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value = $src[$c]
    }
}
